$d = $word.ActiveDocument
$t = $d.Tables.Item(1)
$newValues = @(
    "40+25=",
    "68-37=",
    "31+53=",
    "49+20=",
    "88+4=",
    "99-14=",
    "33+5=",
    "90-30=",
    "9+89=",
    "94-41=",
    "13+0=",
    "21+46=",
    "56-37=",
    "85-71=",
    "72+6=",
    "23+30=",
    "2+49=",
    "78-0=",
    "28+12=",
    "64+30=",
    "24-2=",
    "36+58=",
    "0+60=",
    "87+3=",
    "29+26=",
    "87+11=",
    "97-29=",
    "91-23=",
    "54-47=",
    "66-38=",
    "92-23=",
    "24+18=",
    "13+9=",
    "4+50=",
    "49-11=",
    "65-34=",
    "15+36=",
    "1+58=",
    "17+7=",
    "9-0=",
    "47-39=",
    "0+75=",
    "97-9=",
    "34+56=",
    "20+21=",
    "97-20=",
    "1+88=",
    "32-18=",
    "19+71=",
    "98-82=",
    "56+5=",
    "43-36=",
    "41+30=",
    "37+6=",
    "28-27=",
    "50+37=",
    "73-46=",
    "20+15=",
    "28+16=",
    "54+42=",
    "23+52=",
    "65+24=",
    "13+44=",
    "51+9=",
    "84-56=",
    "3+44=",
    "30+27=",
    "28+70=",
    "59-18=",
    "28+47=",
    "77-53=",
    "13+24=",
    "37+15=",
    "59-1=",
    "84-21=",
    "36+61=",
    "72+0=",
    "86-25=",
    "67-25=",
    "35-27=",
    "78-51=",
    "2+36=",
    "96-26=",
    "77-29=",
    "36-19=",
    "75-73=",
    "36+48=",
    "21+8=",
    "78-59=",
    "73-5=",
    "67-30=",
    "21+11=",
    "75-34=",
    "83-43=",
    "98-39=",
    "45+5=",
    "82-42=",
    "3+50=",
    "83-32=",
    "34+41="
)

$rows = $t.Rows.Count
$cols = $t.Columns.Count
$idx = 0
for ($r = 1; $r -le $rows; $r++) {
    for ($c = 1; $c -le $cols; $c++) {
        $cell = $t.Cell($r, $c)
        $cell.Range.Text = $newValues[$idx]
        $idx = $idx + 1
    }
}
Write-Host "Done. Updated" $idx "cells."